# Insert a new weekly record into the Arveja Verde sheet.
# A new row is inserted at row 70 (pushing the existing rows 70-110 down to
# 71-111), and is populated with a new price observation that mirrors the
# most recent existing record (row 2) but with an updated date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 70; everything below shifts down by one.
$ws.Rows("70:70").Insert()

# Copy formatting (incl. the date number format) from row 2, the most
# recent existing record, into the newly inserted row 70.
$ws.Range("A2:R2").Copy()
$ws.Range("A70").PasteSpecial()

# Populate the new row's values.
$ws.Range("A70").Value = 5
$ws.Range("B70").Value = "Macroferia Regional de Talca"
$ws.Range("C70").Value = "Maule"
$ws.Range("D70").Value = 44875
$ws.Range("E70").Value = 7
$ws.Range("F70").Value = 100112022
$ws.Range("G70").Value = "Arveja Verde"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 500
$ws.Range("K70").Value = 14000
$ws.Range("L70").Value = 14000
$ws.Range("M70").Value = 14000
$ws.Range("N70").Value = "$/saco 25 kilos"
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 560
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"
